$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 134.83333
$ws.Range("I12").Value = 149.8
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = 149.8
$ws.Range("L12").Value = 60
$ws.Range("M12").Value = 20.19999999999999
$ws.Range("N12").Value = -400
$ws.Range("H53").Value = 447.57574
$ws.Range("I53").Value = 425.65518
$ws.Range("K53").Value = 425.65518
$ws.Range("M53").Value = 211.34482
$ws.Range("H76").Value = 12000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 12000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 12000
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -12630
$ws.Range("H79").Value = 12000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 12000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 12000
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -14184
$ws.Range("H106").Value = 20681.428
$ws.Range("I106").Value = 5848.1816
$ws.Range("J106").Value = 36998
$ws.Range("K106").Value = 5848.1816
$ws.Range("L106").Value = 36998
$ws.Range("M106").Value = -5217.1816
$ws.Range("N106").Value = -38260
$ws.Range("H137").Value = 15239
$ws.Range("I137").Value = 7217.9414
$ws.Range("J137").Value = 23761.375
$ws.Range("K137").Value = 21653.8242
$ws.Range("L137").Value = 71284.125
$ws.Range("M137").Value = -19103.8242
$ws.Range("N137").Value = -76384.125
$ws.Range("H138").Value = 4951.362
$ws.Range("I138").Value = 4572.364
$ws.Range("K138").Value = 13717.092
$ws.Range("M138").Value = -8577.091999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 12035.9
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 7632.3706
$ws.Range("I61").Value = 5521.476
$ws.Range("K61").Value = 5521.476
$ws.Range("M61").Value = -5309.476
$ws.Range("H63").Value = 1827.2
$ws.Range("I63").Value = 1111.75
$ws.Range("J63").Value = 4689
$ws.Range("K63").Value = 1111.75
$ws.Range("L63").Value = 4689
$ws.Range("M63").Value = -425.75
$ws.Range("N63").Value = -6061
$ws.Range("H66").Value = 1827.2
$ws.Range("I66").Value = 1111.75
$ws.Range("J66").Value = 4689
$ws.Range("K66").Value = 5558.75
$ws.Range("L66").Value = 23445
$ws.Range("M66").Value = -2126.75
$ws.Range("N66").Value = -30309
$ws.Range("H132").Value = 8431.188
$ws.Range("I132").Value = 5858.2856
$ws.Range("J132").Value = 39949.25
$ws.Range("K132").Value = 17574.8568
$ws.Range("L132").Value = 119847.75
$ws.Range("M132").Value = -15044.8568
$ws.Range("N132").Value = -124907.75
$ws.Range("H136").Value = 7632.3706
$ws.Range("I136").Value = 5521.476
$ws.Range("K136").Value = 16564.428
$ws.Range("M136").Value = -14014.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 27998.625
$ws.Range("J50").Value = 27998.625
$ws.Range("L50").Value = 27998.625
$ws.Range("N50").Value = -29248.625
$ws.Range("H59").Value = 63905.7
$ws.Range("I59").Value = 79998
$ws.Range("J59").Value = 63058.74
$ws.Range("K59").Value = 79998
$ws.Range("L59").Value = 63058.74
$ws.Range("M59").Value = -78853
$ws.Range("N59").Value = -65348.74
$ws.Range("H62").Value = 114032.78
$ws.Range("I62").Value = 169249.17
$ws.Range("J62").Value = 3600
$ws.Range("K62").Value = 169249.17
$ws.Range("L62").Value = 3600
$ws.Range("M62").Value = -168625.17
$ws.Range("N62").Value = -4848
$ws.Range("H65").Value = 114032.78
$ws.Range("I65").Value = 169249.17
$ws.Range("J65").Value = 3600
$ws.Range("K65").Value = 846245.8500000001
$ws.Range("L65").Value = 18000
$ws.Range("M65").Value = -843125.8500000001
$ws.Range("N65").Value = -24240
$ws.Range("H68").Value = 79999
$ws.Range("J68").Value = 79999
$ws.Range("L68").Value = 79999
$ws.Range("N68").Value = -81497
$ws.Range("H71").Value = 79999
$ws.Range("J71").Value = 79999
$ws.Range("L71").Value = 239997
$ws.Range("N71").Value = -247485
$ws.Range("H74").Value = 59999
$ws.Range("J74").Value = 59999
$ws.Range("L74").Value = 59999
$ws.Range("N74").Value = -61747
$ws.Range("H77").Value = 59999
$ws.Range("J77").Value = 59999
$ws.Range("L77").Value = 179997
$ws.Range("N77").Value = -188733

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2355.2
$ws.Range("I68").Value = 1820.3334
$ws.Range("J68").Value = 2848.923
$ws.Range("K68").Value = 5461.0002
$ws.Range("L68").Value = 8546.769
$ws.Range("M68").Value = -4650.0002
$ws.Range("N68").Value = -10168.769
$ws.Range("H71").Value = 2355.2
$ws.Range("I71").Value = 1820.3334
$ws.Range("J71").Value = 2848.923
$ws.Range("K71").Value = 16383.0006
$ws.Range("L71").Value = 25640.307
$ws.Range("M71").Value = -12327.0006
$ws.Range("N71").Value = -33752.307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 15557.556
$ws.Range("I5").Value = 15000
$ws.Range("J5").Value = 15716.857
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 15716.857
$ws.Range("M5").Value = -14888
$ws.Range("N5").Value = -15940.857
$ws.Range("H20").Value = 13498.5
$ws.Range("J20").Value = 13498.5
$ws.Range("L20").Value = 13498.5
$ws.Range("N20").Value = -13988.5
$ws.Range("H132").Value = 2897
$ws.Range("I132").Value = 2866.3333
$ws.Range("K132").Value = 8598.999899999999
$ws.Range("M132").Value = -6068.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 50006
$ws.Range("I23").Value = 50006
$ws.Range("K23").Value = 50006
$ws.Range("M23").Value = -49776
$ws.Range("H25").Value = 10000.2
$ws.Range("I25").Value = 10000.2
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 10000.2
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -9770.200000000001
$ws.Range("N25").ClearContents()
$ws.Range("H68").Value = 2431.457
$ws.Range("I68").Value = 2077.1072
$ws.Range("J68").Value = 3848.8572
$ws.Range("K68").Value = 2077.1072
$ws.Range("L68").Value = 3848.8572
$ws.Range("M68").Value = -1328.1072
$ws.Range("N68").Value = -5346.8572
$ws.Range("H71").Value = 2431.457
$ws.Range("I71").Value = 2077.1072
$ws.Range("J71").Value = 3848.8572
$ws.Range("K71").Value = 10385.536
$ws.Range("L71").Value = 19244.286
$ws.Range("M71").Value = -6641.536
$ws.Range("N71").Value = -26732.286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 6000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 6000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 6000
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -6224
$ws.Range("H21").Value = 23762.5
$ws.Range("J21").Value = 23762.5
$ws.Range("L21").Value = 23762.5
$ws.Range("N21").Value = -24232.5
$ws.Range("J24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("N24").Value = -15460
$ws.Range("H35").Value = 23762.5
$ws.Range("J35").Value = 23762.5
$ws.Range("L35").Value = 23762.5
